$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 2026
$ws1.Range("F11").Value = 37
$ws1.Range("F16").Value = 1388
$ws1.Range("F18").Value = 16
$ws1.Range("F23").Value = 7067
$ws1.Range("F24").Value = 7067
$ws1.Range("F25").Value = 7670
$ws1.Range("F30").Value = 82
$ws1.Range("F32").Value = 252
$ws1.Range("F33").Value = 184
$ws1.Range("F38").Value = 1394
$ws1.Range("F42").Value = 696
$ws1.Range("F46").Value = 224
$ws1.Range("F49").Value = 145

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2585
$ws3.Range("F4").Value = 266

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 266
$ws4.Range("F12").Value = 2026
$ws4.Range("F14").Value = 37
$ws4.Range("F18").Value = 1388
$ws4.Range("F24").Value = 7067
$ws4.Range("F25").Value = 7067
$ws4.Range("F26").Value = 7670
$ws4.Range("F29").Value = 82
$ws4.Range("F30").Value = 252
$ws4.Range("F34").Value = 1394
$ws4.Range("F41").Value = 696
$ws4.Range("F46").Value = 224
$ws4.Range("F48").Value = 145
